$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.320.54'
$ws.Range("E2").Value = '  +2.46%  '
$ws.Range("D3").Value = '1.801.85'
$ws.Range("E3").Value = '  +3.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '340.26'
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4598'
$ws.Range("E7").Value = '  +20.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3812'
$ws.Range("E8").Value = '  +13.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.35'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.154'
$ws.Range("E10").Value = '  +5.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07593'
$ws.Range("E11").Value = '  +6.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.52'
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.355'
$ws.Range("E14").Value = '  +4.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.573'
$ws.Range("E15").Value = '  +8.79%  '
$ws.Range("D16").Value = '1.802.43'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("E17").Value = '  +4.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06723'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.53'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.47'
$ws.Range("E21").Value = '  +5.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.417'
$ws.Range("E22").Value = '  +4.38%  '
$ws.Range("D23").Value = '28.300.42'
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("E24").Value = '  +3.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.432'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.382'
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.65'
$ws.Range("E27").Value = '  +5.36%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.90'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.363'
$ws.Range("E29").Value = '  +5.02%  '
$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '2.009.49'
$ws.Range("E30").Value = '  +3.25%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '132.81'
$ws.Range("E31").Value = '  +2.67%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.244'
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.030'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.855'
$ws.Range("E34").Value = '  +2.13%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.09476'
$ws.Range("E35").Value = '  +9.41%  '
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2300'
$ws.Range("E36").Value = '  +10.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02359'
$ws.Range("E37").Value = '  +4.84%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '12.10'
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06342'
$ws.Range("E39").Value = '  +5.29%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.268'
$ws.Range("E40").Value = '  +3.40%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6616'
$ws.Range("E41").Value = '  +3.07%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.240'
$ws.Range("E42").Value = '  +4.50%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.372'
$ws.Range("E43").Value = '  +5.96%  '
$ws.Range("B44").Value = 'WEMIXTOKEN'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.487'
$ws.Range("E44").Value = '  -2.82%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.14'
$ws.Range("E45").Value = '  +4.37%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.877'
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.6125'
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.25'
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.038'
$ws.Range("E50").Value = '  +3.64%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07163'
$ws.Range("E51").Value = '  +3.44%  '
